$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.533.77'
$ws.Range("E2").Value = '  +0.07%  '

$ws.Range("D3").Value = '3.143.26'
$ws.Range("E3").Value = '  -1.13%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '''573.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.28%  '

$ws.Range("D6").Value = '''164.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.76%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  -5.14%  '

$ws.Range("D9").Value = '3.162.04'
$ws.Range("E9").Value = '  -0.84%  '

$ws.Range("D10").Value = '''0.117'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.52%  '

$ws.Range("E11").Value = '  -2.48%  '

$ws.Range("E12").Value = '  -1.00%  '

$ws.Range("D13").Value = '3.701.91'
$ws.Range("E13").Value = '  -0.97%  '

$ws.Range("E14").Value = '  -1.76%  '

$ws.Range("D15").Value = '64.564.83'
$ws.Range("E15").Value = '  +0.01%  '

$ws.Range("E16").Value = '  -0.93%  '

$ws.Range("D17").Value = '3.153.05'
$ws.Range("E17").Value = '  -1.01%  '

$ws.Range("E18").Value = '  -2.28%  '

$ws.Range("D19").Value = '''407.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.70%  '

$ws.Range("D20").Value = '''5.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.15%  '

$ws.Range("D21").Value = '''12.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.57%  '

$ws.Range("D22").Value = '''7.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.70%  '

$ws.Range("D23").Value = '''0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.05%  '

$ws.Range("D24").Value = '''69.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.74%  '

$ws.Range("E25").Value = '  -2.23%  '

$ws.Range("D26").Value = '''0.196'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.20%  '

$ws.Range("D27").Value = '''0.0000102'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.88%  '

$ws.Range("D28").Value = '''8.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.49%  '

$ws.Range("E29").Value = '  -0.31%  '

$ws.Range("E32").Value = '  -2.33%  '

$ws.Range("D33").Value = '''163.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.03%  '

$ws.Range("D34").Value = '''4.87'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.79%  '

$ws.Range("E36").Value = '  -0.09%  '

$ws.Range("E37").Value = '  -0.39%  '

$ws.Range("E38").Value = '  -0.85%  '

$ws.Range("D39").Value = '2.649.05'
$ws.Range("E39").Value = '  -1.64%  '

$ws.Range("D40").Value = '''23.81'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.02%  '

$ws.Range("E41").Value = '  -3.02%  '

$ws.Range("D42").Value = '''38.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.74%  '

$ws.Range("E43").Value = '  -3.36%  '

$ws.Range("D44").Value = '''0.0613'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.69%  '

$ws.Range("D45").Value = '''5.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.21%  '

$ws.Range("D46").Value = '''290.69'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.44%  '

$ws.Range("D47").Value = '''21.30'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.27%  '

$ws.Range("E48").Value = '  -3.26%  '

$ws.Range("E49").Value = '  -0.15%  '

$ws.Range("D50").Value = '''0.0976'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.67%  '

$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").Value = '''10.49'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.64%  '
